$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2510
$ws.Range("F3").Value = 538
$ws.Range("F5").Value = 277
$ws.Range("F6").Value = 169
$ws.Range("F7").Value = 445
$ws.Range("F8").Value = 1173
$ws.Range("F10").Value = 281
$ws.Range("F11").Value = 107
$ws.Range("F12").Value = 337
$ws.Range("F13").Value = 5403
$ws.Range("F15").Value = 1605
$ws.Range("F16").Value = 3949
$ws.Range("F17").Value = 391
$ws.Range("F20").Value = 4447
$ws.Range("F21").Value = 5848
$ws.Range("F22").Value = 138
$ws.Range("F23").Value = 1009
$ws.Range("F24").Value = 638
$ws.Range("F25").Value = 3623
$ws.Range("F26").Value = 453
$ws.Range("F27").Value = 61
$ws.Range("F28").Value = 176
$ws.Range("F30").Value = 949
$ws.Range("F31").Value = 1331
$ws.Range("F32").Value = 286
$ws.Range("F33").Value = 339
$ws.Range("F34").Value = 1539
$ws.Range("F35").Value = 180
$ws.Range("F36").Value = 1606
$ws.Range("F37").Value = 145
$ws.Range("F38").Value = 1052
$ws.Range("F39").Value = 28
$ws.Range("F40").Value = 1341
$ws.Range("F41").Value = 594
$ws.Range("F42").Value = 84
$ws.Range("F43").Value = 177
$ws.Range("F44").Value = 2716
$ws.Range("F46").Value = 245
$ws.Range("F49").Value = 3841

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1147
$ws.Range("F22").Value = 63

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3603

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3603
$ws.Range("F3").Value = 2510
$ws.Range("F4").Value = 538
$ws.Range("F6").Value = 277
$ws.Range("F7").Value = 1147
$ws.Range("F8").Value = 169
$ws.Range("F9").Value = 445
$ws.Range("F10").Value = 1173
$ws.Range("F12").Value = 281
$ws.Range("F13").Value = 107
$ws.Range("F14").Value = 337
$ws.Range("F15").Value = 5403
$ws.Range("F17").Value = 1605
$ws.Range("F18").Value = 4447
$ws.Range("F19").Value = 5848
$ws.Range("F20").Value = 138
$ws.Range("F21").Value = 1009
$ws.Range("F22").Value = 638
$ws.Range("F23").Value = 3623
$ws.Range("F24").Value = 453
$ws.Range("F25").Value = 61
$ws.Range("F26").Value = 176
$ws.Range("F28").Value = 949
$ws.Range("F29").Value = 1331
$ws.Range("F30").Value = 287
$ws.Range("F31").Value = 340
$ws.Range("F32").Value = 1539
$ws.Range("F33").Value = 180
$ws.Range("F34").Value = 1606
$ws.Range("F36").Value = 1052
$ws.Range("F38").Value = 594
$ws.Range("F41").Value = 84
$ws.Range("F42").Value = 63
$ws.Range("F43").Value = 2716
$ws.Range("F46").Value = 245
$ws.Range("F49").Value = 3841
